$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("connectionMatrix")
$shp = $ws3.Shapes.Item(1)
$full = "Line1" + [char]13 + [char]10 + "Line2"
$shp.TextFrame.Characters().Text = $full
